$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.13083128559433
$ws.Range("C2").Value = 10.47083485335432
$ws.Range("D2").Value = 5.288219612230457
$ws.Range("E2").Value = 11.81046704551659
$ws.Range("F2").Value = 51.50541655708708
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 35.28633954450277
$ws.Range("J2").Value = 10.52997060057699
$ws.Range("K2").Value = 13.82379104765772
$ws.Range("M2").Value = 17.03169587230134
$ws.Range("B3").Value = 12.04086693378412
$ws.Range("C3").Value = 10.37247870353692
$ws.Range("D3").Value = 5.317610910191197
$ws.Range("E3").Value = 11.82745611112423
$ws.Range("F3").Value = 51.23820060089898
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 35.14514794079982
$ws.Range("J3").Value = 10.54178900088667
$ws.Range("K3").Value = 13.77240958628216
$ws.Range("M3").Value = 17.03829306940237
$ws.Range("B4").Value = 11.9896448975629
$ws.Range("C4").Value = 10.31475015156614
$ws.Range("D4").Value = 5.336914150959911
$ws.Range("E4").Value = 11.83991031635697
$ws.Range("F4").Value = 51.08074179322443
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 35.06197438401553
$ws.Range("J4").Value = 10.5501029878294
$ws.Range("K4").Value = 13.74485088046345
$ws.Range("M4").Value = 17.04609730572863
$ws.Range("B5").Value = 11.9698055924706
$ws.Range("C5").Value = 10.29191967988683
$ws.Range("D5").Value = 5.345097126633174
$ws.Range("E5").Value = 11.84549423694146
$ws.Range("F5").Value = 51.01827354326323
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 35.0289788052053
$ws.Range("J5").Value = 10.55375708814004
$ws.Range("K5").Value = 13.73463357201728
$ws.Range("M5").Value = 17.05022212707956
$ws.Range("B6").Value = 11.96657441124179
$ws.Range("C6").Value = 10.28817133365382
$ws.Range("D6").Value = 5.346475054978184
$ws.Range("E6").Value = 11.84645216902842
$ws.Range("F6").Value = 51.00800422687809
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 35.02355454567325
$ws.Range("J6").Value = 10.55437992443639
$ws.Range("K6").Value = 13.73299844441839
$ws.Range("M6").Value = 17.05096410934256
$ws.Range("B7").Value = 11.98937312047298
$ws.Range("C7").Value = 10.3144394085849
$ws.Range("D7").Value = 5.337023225995003
$ws.Range("E7").Value = 11.83998356320606
$ws.Range("F7").Value = 51.07989240767583
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 35.06152573984598
$ws.Range("J7").Value = 10.55015119071936
$ws.Range("K7").Value = 13.74470897235934
$ws.Range("M7").Value = 17.04614910958174
$ws.Range("B8").Value = 12.09899307288955
$ws.Range("C8").Value = 10.43638441402225
$ws.Range("D8").Value = 5.298093289002688
$ws.Range("E8").Value = 11.8159050986402
$ws.Range("F8").Value = 51.41192946782465
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 35.23693478587299
$ws.Range("J8").Value = 10.53382622564336
$ws.Range("K8").Value = 13.80525296986621
$ws.Range("M8").Value = 17.03319206768032
$ws.Range("B9").Value = 12.34466704629655
$ws.Range("C9").Value = 10.69548322837255
$ws.Range("D9").Value = 5.231693322892033
$ws.Range("E9").Value = 11.78473225718034
$ws.Range("F9").Value = 52.11390353426037
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 35.60819839588357
$ws.Range("J9").Value = 10.51019541000404
$ws.Range("K9").Value = 13.95514050277612
$ws.Range("M9").Value = 17.03752013410457
$ws.Range("B10").Value = 12.54220632280939
$ws.Range("C10").Value = 10.89637817595281
$ws.Range("D10").Value = 5.188927412982427
$ws.Range("E10").Value = 11.77159993853605
$ws.Range("F10").Value = 52.6583417649792
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 35.89670185011058
$ws.Range("J10").Value = 10.49793332263322
$ws.Range("K10").Value = 14.08350743057755
$ws.Range("M10").Value = 17.05874291357861
$ws.Range("B11").Value = 12.63538487519399
$ws.Range("C11").Value = 10.98969237793414
$ws.Range("D11").Value = 5.1707699800362
$ws.Range("E11").Value = 11.76774293534267
$ws.Range("F11").Value = 52.91172146767566
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 36.03115638467165
$ws.Range("J11").Value = 10.493459694706
$ws.Range("K11").Value = 14.14568078525634
$ws.Range("M11").Value = 17.07228812916984
$ws.Range("B12").Value = 12.67111206126497
$ws.Range("C12").Value = 11.02527361381442
$ws.Range("D12").Value = 5.164080058426008
$ws.Range("E12").Value = 11.7665862341402
$ws.Range("F12").Value = 53.00844342986905
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 36.08251348016759
$ws.Range("J12").Value = 10.49192420386938
$ws.Range("K12").Value = 14.16975024953943
$ws.Range("M12").Value = 17.0779737161334
$ws.Range("B13").Value = 12.66339844941501
$ws.Range("C13").Value = 11.01760016211003
$ws.Range("D13").Value = 5.165512594104128
$ws.Range("E13").Value = 11.76682184608973
$ws.Range("F13").Value = 52.98757905716502
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 36.07143347448057
$ws.Range("J13").Value = 10.49224785076447
$ws.Range("K13").Value = 14.16454336683786
$ws.Range("M13").Value = 17.07672453716333
$ws.Range("B14").Value = 12.63831550426035
$ws.Range("C14").Value = 10.99261493679102
$ws.Range("D14").Value = 5.170215873781354
$ws.Range("E14").Value = 11.76764168655431
$ws.Range("F14").Value = 52.91966364308111
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 36.03537282073308
$ws.Range("J14").Value = 10.4933301923534
$ws.Range("K14").Value = 14.14765056175489
$ws.Range("M14").Value = 17.07274475655944
$ws.Range("B15").Value = 12.62300806524357
$ws.Range("C15").Value = 10.97734174653687
$ws.Range("D15").Value = 5.173120960874209
$ws.Range("E15").Value = 11.76818341601093
$ws.Range("F15").Value = 52.8781626283544
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 36.01334159654554
$ws.Range("J15").Value = 10.49401380132732
$ws.Range("K15").Value = 14.13737116180926
$ws.Range("M15").Value = 17.07037937090738
$ws.Range("B16").Value = 12.53618070101222
$ws.Range("C16").Value = 10.89031591838773
$ws.Range("D16").Value = 5.190140091116401
$ws.Range("E16").Value = 11.77189455655532
$ws.Range("F16").Value = 52.64189360735286
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 35.88797789581928
$ws.Range("J16").Value = 10.49824788884063
$ws.Range("K16").Value = 14.07951880970725
$ws.Range("M16").Value = 17.0579357337651
$ws.Range("B17").Value = 12.48374019968294
$ws.Range("C17").Value = 10.83739916630706
$ws.Range("D17").Value = 5.200912536596662
$ws.Range("E17").Value = 11.77471306559575
$ws.Range("F17").Value = 52.49837895251324
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 35.81188029302533
$ws.Range("J17").Value = 10.50112807537923
$ws.Range("K17").Value = 14.0449838602991
$ws.Range("M17").Value = 17.05129638646001
$ws.Range("B18").Value = 12.45389213952863
$ws.Range("C18").Value = 10.80714606727995
$ws.Range("D18").Value = 5.207230677029918
$ws.Range("E18").Value = 11.77653349540043
$ws.Range("F18").Value = 52.41637378431307
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 35.76841461912877
$ws.Range("J18").Value = 10.50288865529539
$ws.Range("K18").Value = 14.02547705044087
$ws.Range("M18").Value = 17.04784405743345
$ws.Range("B19").Value = 12.4438411230376
$ws.Range("C19").Value = 10.79693532027981
$ws.Range("D19").Value = 5.209390882605481
$ws.Range("E19").Value = 11.77718410628291
$ws.Range("F19").Value = 52.38870256575622
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 35.7537506358715
$ws.Range("J19").Value = 10.50350262190967
$ws.Range("K19").Value = 14.01893417645617
$ws.Range("M19").Value = 17.04673818551198
$ws.Range("B20").Value = 12.48929029608452
$ws.Range("C20").Value = 10.84301351227538
$ws.Range("D20").Value = 5.199753157378059
$ws.Range("E20").Value = 11.77439240759066
$ws.Range("F20").Value = 52.51360072483248
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 35.8199497276815
$ws.Range("J20").Value = 10.50081071553268
$ws.Range("K20").Value = 14.0486233636142
$ws.Range("M20").Value = 17.05196525205968
$ws.Range("B21").Value = 12.6456712426711
$ws.Range("C21").Value = 10.99994730353185
$ws.Range("D21").Value = 5.168829365902335
$ws.Range("E21").Value = 11.76739263776822
$ws.Range("F21").Value = 52.93959147544117
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 36.04595287350991
$ws.Range("J21").Value = 10.49300798077404
$ws.Range("K21").Value = 14.15259826411013
$ws.Range("M21").Value = 17.07389864434228
$ws.Range("B22").Value = 12.75043794322532
$ws.Range("C22").Value = 11.10392682654414
$ws.Range("D22").Value = 5.14970219881374
$ws.Range("E22").Value = 11.76458868904577
$ws.Range("F22").Value = 53.2224805492042
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 36.19622601019963
$ws.Range("J22").Value = 10.48883260940936
$ws.Range("K22").Value = 14.22360759230098
$ws.Range("M22").Value = 17.09147415940582
$ws.Range("B23").Value = 12.69429906556659
$ws.Range("C23").Value = 11.04831214002078
$ws.Range("D23").Value = 5.159811800934291
$ws.Range("E23").Value = 11.7659233921785
$ws.Range("F23").Value = 53.0711039448442
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 36.11579430544911
$ws.Range("J23").Value = 10.49097660899441
$ws.Range("K23").Value = 14.18543504136603
$ws.Range("M23").Value = 17.08179842145187
$ws.Range("B24").Value = 12.4867801583005
$ws.Range("C24").Value = 10.84047473624942
$ws.Range("D24").Value = 5.200276923578744
$ws.Range("E24").Value = 11.77453675397154
$ws.Range("F24").Value = 52.50671739002519
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 35.81630065054788
$ws.Range("J24").Value = 10.50095386770708
$ws.Range("K24").Value = 14.04697685950993
$ws.Range("M24").Value = 17.05166172162981
$ws.Range("B25").Value = 12.27508518096452
$ws.Range("C25").Value = 10.62341790507823
$ws.Range("D25").Value = 5.248596386053736
$ws.Range("E25").Value = 11.79144835106124
$ws.Range("F25").Value = 51.91878645273502
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 35.50492947291677
$ws.Range("J25").Value = 10.51569179968752
$ws.Range("K25").Value = 13.91132861611723
$ws.Range("M25").Value = 17.03317275678762
